$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 167 (match id 165) - odds refreshed; now shows the fixture that used to be on row 169
$ws.Range("B167").Value = 7630604
$ws.Range("E167").Value = "Treaty United FC"
$ws.Range("F167").Value = "Athlone Town"
$ws.Range("G167").Value = 0
$ws.Range("I167").Value = "A"
$ws.Range("J167").Value = 2.45
$ws.Range("K167").Value = 3.3
$ws.Range("L167").Value = 2.7
$ws.Range("M167").Value = 2.7
$ws.Range("N167").Value = 3.1
$ws.Range("O167").Value = 2.55
$ws.Range("P167").Value = 0
$ws.Range("Q167").Value = 1.95
$ws.Range("R167").Value = 1.85
$ws.Range("S167").Value = 2.25
$ws.Range("T167").Value = 1.825
$ws.Range("U167").Value = 1.975
$ws.Range("W167").Value = -1
$ws.Range("X167").Value = 1.55
$ws.Range("Z167").Value = 0.8500000000000001
$ws.Range("AB167").Value = 0.9750000000000001

# Row 169 (match id 167) - now shows the fixture that used to be on row 167
$ws.Range("B169").Value = 7630603
$ws.Range("E169").Value = "Cork City"
$ws.Range("F169").Value = "Wexford FC"
$ws.Range("G169").Value = 1
$ws.Range("I169").Value = "D"
$ws.Range("J169").Value = 1.6
$ws.Range("K169").Value = 3.75
$ws.Range("L169").Value = 5.25
$ws.Range("M169").Value = 1.666
$ws.Range("N169").Value = 3.75
$ws.Range("O169").Value = 4.75
$ws.Range("P169").Value = -0.75
$ws.Range("Q169").Value = 1.875
$ws.Range("R169").Value = 1.925
$ws.Range("S169").Value = 2.5
$ws.Range("T169").Value = 1.975
$ws.Range("U169").Value = 1.825
$ws.Range("W169").Value = 2.75
$ws.Range("X169").Value = -1
$ws.Range("Z169").Value = 0.925
$ws.Range("AB169").Value = 0.825

# Row 170 (match id 168) - odds refreshed; id now matches what used to be row 172's id
$ws.Range("B170").NumberFormat = "@"
$ws.Range("B170").Value = "7630606"
$ws.Range("B170").Style = "Normal"
$ws.Range("E170").Value = "Treaty United FC"
$ws.Range("F170").Value = "Wexford FC"
$ws.Range("J170").Value = 3.75
$ws.Range("L170").Value = 1.727
$ws.Range("M170").Value = 2.875
$ws.Range("N170").Value = 3.6
$ws.Range("O170").Value = 2.05
$ws.Range("P170").Value = 0.25
$ws.Range("Q170").Value = 1.95
$ws.Range("R170").Value = 1.9
$ws.Range("T170").Value = 1.9
$ws.Range("U170").Value = 1.95

# Row 171 (match id 169) - odds refreshed only
$ws.Range("M171").Value = 4.75
$ws.Range("N171").Value = 3.6
$ws.Range("O171").Value = 1.571
$ws.Range("P171").Value = 0.75
$ws.Range("Q171").Value = 2.05
$ws.Range("R171").Value = 1.8

# Row 172 (match id 170) - odds refreshed; id now matches what used to be row 170's id
$ws.Range("B172").NumberFormat = "@"
$ws.Range("B172").Value = "7630160"
$ws.Range("B172").Style = "Normal"
$ws.Range("E172").Value = "Finn Harps"
$ws.Range("F172").Value = "Longford Town"
$ws.Range("J172").Value = 1.571
$ws.Range("L172").Value = 4.8
$ws.Range("M172").Value = 1.8
$ws.Range("N172").Value = 3.5
$ws.Range("O172").Value = 3.75
$ws.Range("P172").Value = -0.5
$ws.Range("Q172").Value = 1.825
$ws.Range("R172").Value = 2.025
$ws.Range("T172").Value = 2
$ws.Range("U172").Value = 1.85

# Row 173 (match id 171) - odds refreshed only
$ws.Range("T173").Value = 1.8
$ws.Range("U173").Value = 2.05

# Row 174 (match id 172) - odds refreshed only
$ws.Range("T174").Value = 1.925
$ws.Range("U174").Value = 1.925
